{"js": "// Insert the kNN \"result2\" example (prompt + output) right after the\n// existing blank paragraph that follows \"result: Int = 2\", matching the\n// same Normal / justified paragraph formatting used throughout this\n// REPL-transcript section (inherited automatically from the anchor\n// paragraph by insertParagraph).\nconst body = context.document.body;\n\nconst matches = body.search(\"result: Int = 2\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find anchor text \"result: Int = 2\" in the document.');\n}\n\nconst resultPara = matches.items[0].paragraphs.getFirst();\nconst blankPara = resultPara.getNext();\n\nconst promptPara = blankPara.insertParagraph(\n  \"scala> val result2 = classifykNN(Vector(Array(1.0, 0.9)), dataSet, 3)\",\n  Word.InsertLocation.after\n);\n\npromptPara.insertParagraph(\"result2: Int = 1\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert the kNN \"result2\" example (prompt + output) right after the\n# existing blank paragraph that follows \"result: Int = 2\", keeping the\n# same Normal / justified paragraph formatting used throughout this\n# REPL-transcript section.\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"result: Int = 2\")\nif (-not $found) {\n    throw 'Could not find anchor text \"result: Int = 2\" in the document.'\n}\n\n$resultPara = $searchRange.Paragraphs(1)\n$blankPara = $resultPara.Next()\n\n# Add a new paragraph mark after the blank paragraph, then fill it in.\n$insertRange = $blankPara.Range\n$insertRange.Collapse(0)\n$insertRange.InsertParagraphAfter()\n$promptPara = $blankPara.Next()\n$promptPara.Range.InsertAfter(\"scala> val result2 = classifykNN(Vector(Array(1.0, 0.9)), dataSet, 3)\")\n\n# Add a second new paragraph mark after the one we just filled in.\n$promptRange = $promptPara.Range\n$promptRange.Collapse(0)\n$promptRange.InsertParagraphAfter()\n$outputPara = $promptPara.Next()\n$outputPara.Range.InsertAfter(\"result2: Int = 1\")\n"}
